$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.343.15"
$ws.Range("E2").Value = "'  +4.05%  "
$ws.Range("D3").Value = "'2.425.37"
$ws.Range("E3").Value = "'  +0.05%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'318.21"
$ws.Range("E5").Value = "'  +3.82%  "
$ws.Range("D6").Value = "'102.45"
$ws.Range("E6").Value = "'  +5.55%  "
$ws.Range("E7").Value = "'  +1.37%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "'  -0.03%  "
$ws.Range("E9").Value = "'  +8.07%  "
$ws.Range("D10").Value = "'35.58"
$ws.Range("E10").Value = "'  +1.41%  "
$ws.Range("E11").Value = "'  +0.81%  "
$ws.Range("E12").Value = "'  -1.97%  "
$ws.Range("D13").Value = "'18.15"
$ws.Range("E13").Value = "'  -1.77%  "
$ws.Range("D14").Value = "'7.03"
$ws.Range("D15").Value = "'2.805.46"
$ws.Range("E15").Value = "'  +0.37%  "
$ws.Range("D16").Value = "'2.404.15"
$ws.Range("E16").Value = "'  -0.50%  "
$ws.Range("D17").Value = "'0.842"
$ws.Range("E17").Value = "'  +1.87%  "
$ws.Range("D18").Value = "'45.248.00"
$ws.Range("E18").Value = "'  +3.76%  "
$ws.Range("D19").Value = "'12.22"
$ws.Range("E19").Value = "'  +1.48%  "
$ws.Range("D20").Value = "'6.34"
$ws.Range("E20").Value = "'  -1.09%  "
$ws.Range("D21").Value = "'0.0₃0919"
$ws.Range("E22").Value = "'  +1.02%  "
$ws.Range("D23").Value = "'243.89"
$ws.Range("E23").Value = "'  +2.59%  "
$ws.Range("D24").Value = "'2.25"
$ws.Range("E24").Value = "'  -0.09%  "
$ws.Range("E25").Value = "'  +1.97%  "
$ws.Range("E26").Value = "'  -0.06%  "
$ws.Range("D27").Value = "'25.56"
$ws.Range("E27").Value = "'  +2.33%  "
$ws.Range("B28").Value = "'Cosmos"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'9.58"
$ws.Range("E28").Value = "'  +1.57%  "
$ws.Range("B29").Value = "'Toncoin"
$ws.Range("C29").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.07"
$ws.Range("E29").Value = "'  -11.88%  "
$ws.Range("D30").Value = "'49.08"
$ws.Range("D31").Value = "'32.87"
$ws.Range("E31").Value = "'  +1.71%  "
$ws.Range("E32").Value = "'  +5.94%  "
$ws.Range("D33").Value = "'20.18"
$ws.Range("E33").Value = "'  +9.36%  "
$ws.Range("E34").Value = "'  +1.46%  "
$ws.Range("D36").Value = "'0.0765"
$ws.Range("E36").Value = "'  +1.91%  "
$ws.Range("D37").Value = "'1.87"
$ws.Range("E37").Value = "'  -1.47%  "
$ws.Range("E38").Value = "'  +1.61%  "
$ws.Range("B39").Value = "'Monero"
$ws.Range("C39").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'126.89"
$ws.Range("E39").Value = "'  -2.37%  "
$ws.Range("B40").Value = "'LidoDAOToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.86"
$ws.Range("E40").Value = "'  -2.46%  "
$ws.Range("E41").Value = "'  -2.26%  "
$ws.Range("E42").Value = "'  +0.87%  "
$ws.Range("D43").Value = "'20.50"
$ws.Range("E43").Value = "'  -2.90%  "
$ws.Range("E44").Value = "'  +2.44%  "
$ws.Range("D45").Value = "'1.933.89"
$ws.Range("E45").Value = "'  -0.53%  "
$ws.Range("E46").Value = "'  -2.94%  "
$ws.Range("E47").Value = "'  +3.33%  "
$ws.Range("D48").Value = "'1.82"
$ws.Range("E48").Value = "'  +15.94%  "
$ws.Range("E49").Value = "'  -2.22%  "
$ws.Range("D50").Value = "'76.41"
$ws.Range("E50").Value = "'  +5.56%  "
$ws.Range("D51").Value = "'53.92"
$ws.Range("E51").Value = "'  +2.46%  "
